# "adding occupancy type to data helper"
#
# Insert a new "Occ_m2pax" column right after the Code column on the
# INTERNAL_LOADS sheet, populate it, and make INTERNAL_LOADS the active
# sheet/tab (INDOOR_COMFORT loses that status).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. INTERNAL_LOADS: insert new column B ("Occ_m2pax") with its values
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("INTERNAL_LOADS")

# Insert a new column before the current column B - everything from the
# old B (X_ghp) onward shifts one column to the right.
$ws.Columns.Item(2).Insert()

# The freshly inserted column comes back formatted as text (inherited from
# column A's "@" cell style), which would turn any numbers typed into it
# into shared strings. Re-base its number formatting on column I (the old
# "Vww_lpd" column, now shifted one to the right), which is consistently
# formatted as a plain number column across every data row, *before*
# writing any values.
$ws.Range("I2:I20").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)  # xlPasteFormats

# Header
$ws.Range("B1").Value = "Occ_m2pax"

# Data values (rows 2-20)
$occValues = @(35, 60, 23, 10, 6, 0, 2.7, 13, 19, 4, 19, 9, 20, 0, 0, 0, 20, 10, 9)
for ($i = 0; $i -lt $occValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $occValues[$i]
}

# Column width to match the "Code" column (A) - same style, no bestFit
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Selection / active view bookkeeping
$ws.Range("B21").Select()
$ws.Activate()

# ---------------------------------------------------------------------
# 2. INDOOR_COMFORT loses the "active sheet" flag
# ---------------------------------------------------------------------
$wsIndoor = $wb.Worksheets.Item("INDOOR_COMFORT")
$wsIndoor.Range("Q7").Select()

$ws.Activate()
